# Apply updated crypto price/volume data (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force the cell to store a literal text value (matches the source
    # inlineStr cells) instead of letting Excel auto-coerce numeric-
    # looking strings ("0.9988", "1.000", ...) into real numbers, then
    # drop back to the default style so no formatting diff is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.219.97"
$ws.Range("E2").Value = "  +0.37%  "

Set-TextValue $ws.Range("D3") "1.832.87"
$ws.Range("E3").Value = "  -0.16%  "

Set-TextValue $ws.Range("D4") "0.9988"
$ws.Range("E4").Value = "  -0.26%  "

Set-TextValue $ws.Range("D5") "243.06"
$ws.Range("E5").Value = "  -0.30%  "

Set-TextValue $ws.Range("D6") "0.6194"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("E7").Value = "  -0.34%  "

Set-TextValue $ws.Range("D10") "23.30"
$ws.Range("E10").Value = "  +0.81%  "

Set-TextValue $ws.Range("D11") "0.07671"
$ws.Range("E11").Value = "  -0.38%  "

Set-TextValue $ws.Range("D12") "1.832.14"
$ws.Range("E12").Value = "  -0.33%  "

Set-TextValue $ws.Range("D13") "4.981"
$ws.Range("E13").Value = "  -0.17%  "

Set-TextValue $ws.Range("D14") "0.6702"
$ws.Range("E14").Value = "  -0.08%  "

Set-TextValue $ws.Range("D15") "82.69"
$ws.Range("E15").Value = "  +0.15%  "

Set-TextValue $ws.Range("D16") "0.000008954"
$ws.Range("E16").Value = "  -3.88%  "

Set-TextValue $ws.Range("D17") "5.876"
$ws.Range("E17").Value = "  -0.91%  "

Set-TextValue $ws.Range("D18") "29.199.28"
$ws.Range("E18").Value = "  +0.35%  "

Set-TextValue $ws.Range("D19") "2.072.77"
$ws.Range("E19").Value = "  -1.18%  "

Set-TextValue $ws.Range("D20") "236.03"
$ws.Range("E20").Value = "  +2.15%  "

Set-TextValue $ws.Range("D21") "12.51"
$ws.Range("E21").Value = "  -0.92%  "

Set-TextValue $ws.Range("D22") "0.9999"
$ws.Range("E22").Value = "  -0.39%  "

Set-TextValue $ws.Range("D23") "7.347"
$ws.Range("E23").Value = "  +2.58%  "

Set-TextValue $ws.Range("D24") "0.9996"
$ws.Range("E24").Value = "  -0.36%  "

Set-TextValue $ws.Range("D25") "158.18"
$ws.Range("E25").Value = "  -1.22%  "

Set-TextValue $ws.Range("D26") "0.1399"
$ws.Range("E26").Value = "  +0.97%  "

Set-TextValue $ws.Range("D27") "8.558"
$ws.Range("E27").Value = "  +0.59%  "

Set-TextValue $ws.Range("D28") "17.63"
$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("E29").Value = "  -0.76%  "

Set-TextValue $ws.Range("D30") "0.05764"
$ws.Range("E30").Value = "  +4.36%  "

Set-TextValue $ws.Range("D31") "4.112"
$ws.Range("E31").Value = "  -0.24%  "

$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("E33").Value = "  +0.72%  "

Set-TextValue $ws.Range("D34") "1.879"
$ws.Range("E34").Value = "  +2.38%  "

Set-TextValue $ws.Range("D35") "0.7328"
$ws.Range("E35").Value = "  -1.61%  "

Set-TextValue $ws.Range("D36") "1.142"
$ws.Range("E36").Value = "  +0.23%  "

Set-TextValue $ws.Range("D37") "2.609"
$ws.Range("E37").Value = "  -2.08%  "

Set-TextValue $ws.Range("D38") "2.860"
$ws.Range("E38").Value = "  +3.32%  "

Set-TextValue $ws.Range("D39") "1.224.28"
$ws.Range("E39").Value = "  +0.42%  "

Set-TextValue $ws.Range("D40") "0.01752"
$ws.Range("E40").Value = "  -1.56%  "

Set-TextValue $ws.Range("D41") "6.255"
$ws.Range("E41").Value = "  -2.95%  "

Set-TextValue $ws.Range("D42") "0.9088"
$ws.Range("E42").Value = "  +1.55%  "

Set-TextValue $ws.Range("D43") "1.000"
$ws.Range("E43").Value = "  -0.19%  "

Set-TextValue $ws.Range("D44") "101.78"
$ws.Range("E44").Value = "  -0.04%  "

Set-TextValue $ws.Range("D45") "1.977.56"
$ws.Range("E45").Value = "  -1.20%  "

Set-TextValue $ws.Range("D46") "65.60"
$ws.Range("E46").Value = "  +0.07%  "

Set-TextValue $ws.Range("D47") "0.5037"
$ws.Range("E47").Value = "  -1.22%  "

Set-TextValue $ws.Range("D48") "0.00000000119"
$ws.Range("E48").Value = "  -2.51%  "

Set-TextValue $ws.Range("D49") "9.126"
$ws.Range("E49").Value = "  +0.26%  "

Set-TextValue $ws.Range("D50") "0.4024"
$ws.Range("E50").Value = "  -1.11%  "

Set-TextValue $ws.Range("D51") "0.1135"
$ws.Range("E51").Value = "  +3.27%  "

# Rows 8 and 9: Cardano and Dogecoin swapped ranking positions
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
Set-TextValue $ws.Range("D8") "0.2947"
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D9") "0.07370"
$ws.Range("E9").Value = "  -1.29%  "
